$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C to hold Industry data
$ws.Columns.Item(3).Insert()

# Set header for new Industry column
$ws.Range("C1").Value = "Industry"

# Populate Industry values for each holding row
$ws.Range("C2").Value = "Retailing"
$ws.Range("C3").Value = "Automobiles"
$ws.Range("C4").Value = "Automobiles"
$ws.Range("C5").Value = "Telecom - Services"
$ws.Range("C6").Value = "Transport Services"
$ws.Range("C7").Value = "Auto Components"
$ws.Range("C8").Value = "Retailing"
$ws.Range("C9").Value = "Healthcare Services"
$ws.Range("C10").Value = "Transport Services"
$ws.Range("C11").Value = "Consumer Durables"
$ws.Range("C12").Value = "Personal Products"
$ws.Range("C13").Value = "Industrial Manufacturing"
$ws.Range("C14").Value = "Healthcare Services"
$ws.Range("C15").Value = "Household Products"
$ws.Range("C16").Value = "Agricultural Food & other Products"
$ws.Range("C17").Value = "Consumer Durables"
$ws.Range("C18").Value = "Retailing"
$ws.Range("C19").Value = "Retailing"
$ws.Range("C20").Value = "Leisure Services"
$ws.Range("C21").Value = "Financial Technology (Fintech)"
$ws.Range("C22").Value = "Healthcare Services"
$ws.Range("C23").Value = "Realty"
$ws.Range("C24").Value = "Automobiles"
$ws.Range("C25").Value = "Retailing"
$ws.Range("C26").Value = "Consumer Durables"
$ws.Range("C27").Value = "Consumer Durables"
$ws.Range("C28").Value = "Auto Components"
$ws.Range("C29").Value = "Pharmaceuticals & Biotechnology"
$ws.Range("C30").Value = "Insurance"
$ws.Range("C31").Value = "Commercial Services & Supplies"
$ws.Range("C32").Value = "Consumer Durables"
$ws.Range("C33").Value = "Industrial Products"
$ws.Range("C34").Value = "Cigarettes & Tobacco Products"
$ws.Range("C35").Value = "Automobiles"
$ws.Range("C36").Value = "Retailing"
